$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values updated to reflect a newer handback report generation run.

# Shared value (same string used on both Overview and de-de sheets)
$wsOverview.Range("G2").Value = "2016-08-25 11:08:04"
$wsDeDe.Range("H2").Value     = "2016-08-25 11:08:04"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-25 11:07:55"
$wsZhCn.Range("K2").Value = "2016-08-25 11:08:27"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-25 11:08:34"
